$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of new Price (D) values are plain decimals (e.g. '529.17',
# '1.00', '0.150') that Excel's Range.Value setter would otherwise parse
# as numbers, silently dropping trailing zeros. Force just those specific
# cells to Text format first so the literal string is preserved exactly,
# matching the original inline-string cells; every other cell (including
# the two-part-decimal and Volume(1h) values, which are never ambiguous)
# is left with its original formatting untouched.
$textFormatCells = @("D5", "D6", "D7", "D8", "D9", "D12", "D13", "D16", "D17", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D33", "D34", "D36", "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49")
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "71.103.68"
$ws.Range("E2").Value = "  +2.48%  "

# Row 3
$ws.Range("D3").Value = "4.003.18"
$ws.Range("E3").Value = "  +1.40%  "

# Row 4
$ws.Range("E4").Value = "  +0.18%  "

# Row 5
$ws.Range("D5").Value = "529.17"
$ws.Range("E5").Value = "  +5.07%  "

# Row 6
$ws.Range("D6").Value = "147.85"
$ws.Range("E6").Value = "  -0.36%  "

# Row 7
$ws.Range("D7").Value = "0.621"
$ws.Range("E7").Value = "  -0.98%  "

# Row 8
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.17%  "

# Row 9
$ws.Range("D9").Value = "0.738"
$ws.Range("E9").Value = "  +0.32%  "

# Row 10
$ws.Range("E10").Value = "  +0.10%  "

# Row 11
$ws.Range("E11").Value = "  -0.32%  "

# Row 12
$ws.Range("D12").Value = "42.94"
$ws.Range("E12").Value = "  -1.77%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "10.53"
$ws.Range("E13").Value = "  +0.00%  "

# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "4.608.78"
$ws.Range("E14").Value = "  +0.78%  "

# Row 15
$ws.Range("D15").Value = "4.000.71"
$ws.Range("E15").Value = "  +1.30%  "

# Row 16
$ws.Range("D16").Value = "21.37"
$ws.Range("E16").Value = "  +6.67%  "

# Row 17
$ws.Range("D17").Value = "14.23"
$ws.Range("E17").Value = "  -0.75%  "

# Row 18
$ws.Range("E18").Value = "  +2.18%  "

# Row 19
$ws.Range("E19").Value = "  -1.91%  "

# Row 20
$ws.Range("D20").Value = "71.112.48"
$ws.Range("E20").Value = "  +2.42%  "

# Row 21
$ws.Range("D21").Value = "440.29"
$ws.Range("E21").Value = "  +0.42%  "

# Row 22
$ws.Range("D22").Value = "3.54"
$ws.Range("E22").Value = "  +2.33%  "

# Row 23
$ws.Range("D23").Value = "90.43"
$ws.Range("E23").Value = "  +1.57%  "

# Row 24
$ws.Range("B24").Value = "RenderToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D24").Value = "12.49"
$ws.Range("E24").Value = "  +3.28%  "

# Row 25
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "14.34"
$ws.Range("E25").Value = "  -2.21%  "

# Row 26
$ws.Range("D26").Value = "4.07"
$ws.Range("E26").Value = "  +5.14%  "

# Row 27
$ws.Range("D27").Value = "10.77"
$ws.Range("E27").Value = "  -4.08%  "

# Row 28
$ws.Range("D28").Value = "37.00"
$ws.Range("E28").Value = "  -0.48%  "

# Row 29
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "13.47"
$ws.Range("E29").Value = "  -0.10%  "

# Row 30
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").Value = "692.81"
$ws.Range("E30").Value = "  -1.29%  "

# Row 31
$ws.Range("D31").Value = "0.128"
$ws.Range("E31").Value = "  -1.24%  "

# Row 32
$ws.Range("E32").Value = "  -0.35%  "

# Row 33
$ws.Range("D33").Value = "6.78"
$ws.Range("E33").Value = "  +11.98%  "

# Row 34
$ws.Range("D34").Value = "67.03"
$ws.Range("E34").Value = "  +5.16%  "

# Row 35
$ws.Range("D35").Value = "0.0$([char]0x2083)0928"
$ws.Range("E35").Value = "  +3.75%  "

# Row 36
$ws.Range("D36").Value = "0.444"
$ws.Range("E36").Value = "  -1.78%  "

# Row 37
$ws.Range("D37").Value = "40.29"
$ws.Range("E37").Value = "  -2.17%  "

# Row 38
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.150"
$ws.Range("E38").Value = "  -0.87%  "

# Row 39
$ws.Range("B39").Value = "ThetaToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D39").Value = "3.42"
$ws.Range("E39").Value = "  +10.82%  "

# Row 40
$ws.Range("E40").Value = "  +0.05%  "

# Row 41
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  -0.28%  "

# Row 42
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").Value = "2.94"
$ws.Range("E42").Value = "  +1.12%  "

# Row 43
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.0485"
$ws.Range("E43").Value = "  -0.92%  "

# Row 44
$ws.Range("D44").Value = "3.14"
$ws.Range("E44").Value = "  +2.97%  "

# Row 45
$ws.Range("D45").Value = "3.54"
$ws.Range("E45").Value = "  +4.83%  "

# Row 46
$ws.Range("D46").Value = "3.23"
$ws.Range("E46").Value = "  +7.45%  "

# Row 47
$ws.Range("D47").Value = "0.144"
$ws.Range("E47").Value = "  -0.20%  "

# Row 48
$ws.Range("D48").Value = "0.000283"
$ws.Range("E48").Value = "  +19.37%  "

# Row 49
$ws.Range("D49").Value = "9.23"
$ws.Range("E49").Value = "  +5.17%  "

# Row 50
$ws.Range("E50").Value = "  -0.38%  "

# Row 51
$ws.Range("E51").Value = "  -0.25%  "

